# Priorización.xlsx — remove the two blank header rows above the CUS
# (Caso de Uso) list on the "CUS" sheet, then leave that sheet active
# with F3 selected (matching the saved view state of the edited file).

$wb = $excel.ActiveWorkbook

# --- CUS sheet: delete rows 2:3 (blank rows above the "ID / Caso de Uso" table) ---
$cus = $wb.Worksheets.Item("CUS")
$cus.Range("A2:A3").EntireRow.Delete()

# --- Priorización sheet: reset its remembered selection back to A1 ---
$pri = $wb.Worksheets.Item("Priorización")
$pri.Activate()
$pri.Range("A1").Select()

# --- CUS becomes the active sheet again, with F3 selected ---
$cus.Activate()
$cus.Range("F3").Select()
